# Reservation sheet fix:
#  - Booking status in row 7 (Reservation ID 6) is actually "Cancelled", not "confirmed"
#  - A new customer reservation (row 15 / Reservation ID 14) was missing from the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the booking status that was recorded incorrectly for reservation #6 (row 7)
$ws.Range("H7").Value = "Cancelled"

# Append the new reservation as row 15
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 45413
$ws.Range("C15").NumberFormat = "yyyy-MM-dd"
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = "2 seat"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = "10:00"
$ws.Range("H15").Value = "pending"
